## Actualización automática 2025-10-20 13:30:09
## Applies updated sales figures for GUERRERO FAREZ FABIAN MAURICIO across
## the three report sheets (VENTAS POR GRUPO, VENTA MENSUAL,
## CUMPLIMIENTO MENSUAL), keeping the cross-sheet totals consistent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (per-client / per-group breakdown)
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M11").Value = 2370.29

$wsGrupo.Range("K24").Value = 786.6
$wsGrupo.Range("M24").Value = 3553.92

$wsGrupo.Range("L29").Value = 1882.38

$wsGrupo.Range("K37").Value = 624.24

$wsGrupo.Range("K42").Value = 812.16

$wsGrupo.Range("K56").Value = "6 de 54"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" (monthly totals; "octubre" column updated)
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F11").Value = 2370.29
$wsMensual.Range("F24").Value = 4340.52
$wsMensual.Range("F29").Value = 10142.02
$wsMensual.Range("F37").Value = 3153.57
$wsMensual.Range("F42").Value = 898.5599999999999
$wsMensual.Range("F60").Value = 49066.76

# Column F widened slightly to fit the new, longer numbers.
$wsMensual.Columns.Item(6).ColumnWidth = 13.166666666666666

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (compliance summary, recalculated)
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PANELES DECORATIVOS
$wsCumpl.Range("D10").Value = 3973.97
$wsCumpl.Range("E10").Value = -92.89016465607983
$wsCumpl.Range("F10").Value = 1.023934103032397

# PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 6244.13
$wsCumpl.Range("E11").Value = 5586.87
$wsCumpl.Range("F11").Value = 0.527777026455921

# PORCELANATO
$wsCumpl.Range("D12").Value = 25175.43
$wsCumpl.Range("E12").Value = 27487.69
$wsCumpl.Range("F12").Value = 0.4780466861819049

# TOTAL
$wsCumpl.Range("D14").Value = 47046.07
$wsCumpl.Range("E14").Value = 51970.43661190613
$wsCumpl.Range("F14").Value = 0.4751336076155104
